# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape, leaving all other data untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11588
    3  = 11122
    6  = 1006
    7  = 119
    9  = 41
    10 = 43
    11 = 10673
    12 = 4124
    15 = 2460
    16 = 48
    17 = 42
    18 = 120
    19 = 432
    20 = 11112
    21 = 10879
    23 = 24
    26 = 26
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
